$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.223.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.827.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5973'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.79%  '
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06957'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2745'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.19'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07623'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.839.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.750'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6232'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009698'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '78.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '28.908.09'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.713'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -8.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '221.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("E21").Value = '  -6.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.858'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.06%  '
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '155.87'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.924'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1287'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06712'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.443'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.08%  '
$ws.Range("E30").Value = '  -2.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.825'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.747'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.089'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.710'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6400'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.546'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.736'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.185.91'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01729'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.503'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9023'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.01%  '
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.979.36'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.31'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.88'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000113'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.457'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4556'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05505'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.567'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.326'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.15%  '
